$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 749
$ws.Range("J70").Value = 749
$ws.Range("L70").Value = 2247
$ws.Range("N70").Value = -2787

$ws.Range("H73").Value = 749
$ws.Range("J73").Value = 749
$ws.Range("L73").Value = 2247
$ws.Range("N73").Value = -4119

$ws.Range("H80").Value = 1412.2927
$ws.Range("J80").Value = 2311.3
$ws.Range("L80").Value = 6933.900000000001
$ws.Range("N80").Value = -8929.900000000001

$ws.Range("H83").Value = 1412.2927
$ws.Range("J83").Value = 2311.3
$ws.Range("L83").Value = 20801.7
$ws.Range("N83").Value = -30785.7

$ws.Range("H96").Value = 568.8333
$ws.Range("I96").Value = 536.75
$ws.Range("J96").Value = 633
$ws.Range("K96").Value = 1610.25
$ws.Range("L96").Value = 1899
$ws.Range("M96").Value = -237.25
$ws.Range("N96").Value = -4645

$ws.Range("H113").Value = 8544.809999999999
$ws.Range("I113").Value = 7453.8335
$ws.Range("K113").Value = 7453.8335
$ws.Range("M113").Value = -4199.8335

$ws.Range("H125").Value = 3449.6667
$ws.Range("I125").Value = 2825
$ws.Range("J125").Value = 3628.1428
$ws.Range("K125").Value = 25425
$ws.Range("L125").Value = 32653.2852
$ws.Range("M125").Value = -22965
$ws.Range("N125").Value = -37573.2852

$ws.Range("H135").Value = 13895009
$ws.Range("I135").Value = 25002366
$ws.Range("J135").Value = 10812.3125
$ws.Range("K135").Value = 225021294
$ws.Range("L135").Value = 97310.8125
$ws.Range("M135").Value = -225018759
$ws.Range("N135").Value = -102380.8125

$ws.Range("H138").Value = 5766.6177
$ws.Range("J138").Value = 6622.88
$ws.Range("L138").Value = 19868.64
$ws.Range("N138").Value = -30148.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21260610
$ws.Range("I32").Value = 21758170
$ws.Range("K32").Value = 21758170
$ws.Range("M32").Value = -21757883

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H88").Value = 7577805
$ws.Range("I88").Value = 1270
$ws.Range("K88").Value = 1270
$ws.Range("M88").Value = -864

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H91").Value = 7577805
$ws.Range("I91").Value = 1270
$ws.Range("K91").Value = 1270
$ws.Range("M91").Value = 134

$ws.Range("H122").Value = 3319.28
$ws.Range("I122").Value = 2564.1177
$ws.Range("J122").Value = 4924
$ws.Range("K122").Value = 7692.353099999999
$ws.Range("L122").Value = 14772
$ws.Range("M122").Value = -5242.353099999999
$ws.Range("N122").Value = -19672

$ws.Range("H132").Value = 4963.066
$ws.Range("I132").Value = 4344.311
$ws.Range("J132").Value = 5861.2583
$ws.Range("K132").Value = 13032.933
$ws.Range("L132").Value = 17583.7749
$ws.Range("M132").Value = -10502.933
$ws.Range("N132").Value = -22643.7749

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 1000
$ws.Range("J15").Value = 1000
$ws.Range("L15").Value = 1000
$ws.Range("N15").Value = -1454

$ws.Range("H86").Value = 6127.636
$ws.Range("I86").Value = 6125.5
$ws.Range("K86").Value = 6125.5
$ws.Range("M86").Value = -5002.5

$ws.Range("H89").Value = 6127.636
$ws.Range("I89").Value = 6125.5
$ws.Range("K89").Value = 30627.5
$ws.Range("M89").Value = -25011.5

$ws.Range("H105").Value = 3848.4666
$ws.Range("I105").Value = 3312.4
$ws.Range("J105").Value = 4116.5
$ws.Range("K105").Value = 3312.4
$ws.Range("L105").Value = 4116.5
$ws.Range("M105").Value = -1565.4
$ws.Range("N105").Value = -7610.5

$ws.Range("H107").Value = 6799
$ws.Range("I107").Value = 8027.2856
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 8027.2856
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = -6107.2856
$ws.Range("N107").Value = -6340

$ws.Range("H134").Value = 3694.1526
$ws.Range("I134").Value = 2951.8572
$ws.Range("K134").Value = 8855.571599999999
$ws.Range("M134").Value = -6320.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 966.6667
$ws.Range("I22").Value = 966.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 966.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -616.6667
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 5726.174
$ws.Range("I31").Value = 4633.467
$ws.Range("J31").Value = 7775
$ws.Range("K31").Value = 4633.467
$ws.Range("L31").Value = 7775
$ws.Range("M31").Value = -4338.467
$ws.Range("N31").Value = -8365

$ws.Range("H34").Value = 5726.174
$ws.Range("I34").Value = 4633.467
$ws.Range("J34").Value = 7775
$ws.Range("K34").Value = 4633.467
$ws.Range("L34").Value = 7775
$ws.Range("M34").Value = -4431.467
$ws.Range("N34").Value = -8179

$ws.Range("H35").Value = 295.5238
$ws.Range("I35").Value = 315.16666
$ws.Range("J35").Value = 177.66667
$ws.Range("K35").Value = 315.16666
$ws.Range("L35").Value = 177.66667
$ws.Range("M35").Value = -21.16665999999998
$ws.Range("N35").Value = -765.6666700000001

$ws.Range("H105").Value = 1023.9167
$ws.Range("I105").Value = 935.875
$ws.Range("K105").Value = 935.875
$ws.Range("M105").Value = 811.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 418.875
$ws.Range("I26").Value = 478.42856
$ws.Range("K26").Value = 1435.28568
$ws.Range("M26").Value = -1147.28568

$ws.Range("H68").Value = 1142
$ws.Range("J68").Value = 1217.5
$ws.Range("L68").Value = 3652.5
$ws.Range("N68").Value = -5274.5

$ws.Range("H71").Value = 1142
$ws.Range("J71").Value = 1217.5
$ws.Range("L71").Value = 10957.5
$ws.Range("N71").Value = -19069.5

$ws.Range("H99").Value = 33486
$ws.Range("I99").Value = 1360
$ws.Range("J99").Value = 65612
$ws.Range("K99").Value = 4080
$ws.Range("L99").Value = 196836
$ws.Range("M99").Value = -1834
$ws.Range("N99").Value = -201328

$ws.Range("H108").Value = 1750
$ws.Range("I108").Value = 2500
$ws.Range("J108").Value = 1000
$ws.Range("K108").Value = 7500
$ws.Range("L108").Value = 3000
$ws.Range("M108").Value = -4620
$ws.Range("N108").Value = -8760

$ws.Range("H109").Value = 2145.2
$ws.Range("I109").Value = 581
$ws.Range("K109").Value = 1743
$ws.Range("M109").Value = -703

$ws.Range("H121").Value = 3521.8
$ws.Range("J121").Value = 4120.2856
$ws.Range("L121").Value = 12360.8568
$ws.Range("N121").Value = -14980.8568

$ws.Range("H122").Value = 2029.6111
$ws.Range("I122").Value = 1439.4
$ws.Range("J122").Value = 2256.6155
$ws.Range("K122").Value = 12954.6
$ws.Range("L122").Value = 20309.5395
$ws.Range("M122").Value = -10504.6
$ws.Range("N122").Value = -25209.5395

$ws.Range("H134").Value = 7611.905
$ws.Range("I134").Value = 2989.2856
$ws.Range("K134").Value = 8967.856800000001
$ws.Range("M134").Value = -3897.856800000001

$ws.Range("H139").Value = 31261542
$ws.Range("I139").Value = 90914480
$ws.Range("K139").Value = 272743440
$ws.Range("M139").Value = -272738300

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2431.913
$ws.Range("I102").Value = 2428.818
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2428.818
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -806.8180000000002
$ws.Range("N102").Value = -5744

$ws.Range("H126").Value = 5312.0835
$ws.Range("I126").Value = 3650
$ws.Range("J126").Value = 5644.5
$ws.Range("K126").Value = 10950
$ws.Range("L126").Value = 16933.5
$ws.Range("M126").Value = -8480
$ws.Range("N126").Value = -21873.5

$ws.Range("H132").Value = 12115.19
$ws.Range("I132").Value = 10884.571
$ws.Range("K132").Value = 32653.713
$ws.Range("M132").Value = -30123.713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1499.6666
$ws.Range("I22").Value = 1499.6666
$ws.Range("K22").Value = 1499.6666
$ws.Range("M22").Value = -1204.6666

$ws.Range("H27").Value = 1499.6666
$ws.Range("I27").Value = 1499.6666
$ws.Range("K27").Value = 1499.6666
$ws.Range("M27").Value = -1392.6666

$ws.Range("H46").Value = 8772.532999999999
$ws.Range("I46").Value = 6757.6
$ws.Range("J46").Value = 9175.52
$ws.Range("K46").Value = 6757.6
$ws.Range("L46").Value = 9175.52
$ws.Range("M46").Value = -6569.6
$ws.Range("N46").Value = -9551.52

$ws.Range("H55").Value = 1961.4166
$ws.Range("I55").Value = 2753
$ws.Range("J55").Value = 378.25
$ws.Range("K55").Value = 2753
$ws.Range("L55").Value = 378.25
$ws.Range("M55").Value = -2580
$ws.Range("N55").Value = -724.25

$ws.Range("H132").Value = 7043.6
$ws.Range("I132").Value = 6406.4375
$ws.Range("K132").Value = 19219.3125
$ws.Range("M132").Value = -16689.3125

$ws.Range("H136").Value = 6229.643
$ws.Range("I136").Value = 4142.857
$ws.Range("K136").Value = 12428.571
$ws.Range("M136").Value = -9878.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9666.166999999999
$ws.Range("I122").Value = 11332.667
$ws.Range("J122").Value = 7999.6665
$ws.Range("K122").Value = 33998.001
$ws.Range("L122").Value = 23998.9995
$ws.Range("M122").Value = -31548.001
$ws.Range("N122").Value = -28898.9995

$ws.Range("H132").Value = 5861.528
$ws.Range("I132").Value = 4651
$ws.Range("J132").Value = 6944.6313
$ws.Range("K132").Value = 13953
$ws.Range("L132").Value = 20833.8939
$ws.Range("M132").Value = -11423
$ws.Range("N132").Value = -25893.8939

$ws.Range("H136").Value = 5209.3105
$ws.Range("I136").Value = 6689.4
$ws.Range("J136").Value = 4430.316
$ws.Range("K136").Value = 20068.2
$ws.Range("L136").Value = 13290.948
$ws.Range("M136").Value = -17518.2
$ws.Range("N136").Value = -18390.948
